$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data table: add the "Cocktail Sort" column (N) and rename
#     "Booble sort" -> "bubble sort" (K column header).
# Order matters for shared-string table layout: write the new unique
# string (Cocktail Sort) before renaming the existing one so the
# renamed string lands after it, matching the target shared strings order.
$ws.Range("N3").Value = "Cocktail Sort"
$ws.Range("K3").Value = "bubble sort"

$ws.Range("N4").Value = 0.009
$ws.Range("N5").Value = 0.0703
$ws.Range("N6").Value = 0.22

# --- Sheet view changes: zoom + new selected cell.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("K4").Select()

# --- Chart updates ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# Add the new "Cocktail Sort" series to the line chart.
$newSeries = $chart.SeriesCollection().NewSeries()
$newSeries.Name = '=Hoja1!$N$3'
$newSeries.Values = '=Hoja1!$N$4:$N$6'
$newSeries.XValues = '=Hoja1!$J$4:$J$6'

Write-Output "done"
